$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "audioFalse" header to "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Both rows' "audioFalse" data values now reference the shared "train2P2" string
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
